# bq25570SolarAppDesignExample_V1p3.xlsx - "Computation" sheet value updates.
#
# The underlying XML diff only changes four raw input cells; every other
# changed <c> in the diff is a formula cell (H7, I7, J7, H11, H13, I13, J13,
# I17, I19, I23, J25, J27, J28, J29, J31, J34, J37, ...) whose cached <v>
# simply reflects automatic recalculation once the inputs below change.
# Excel (and this COM host) recalculates those for us, so we only need to
# touch the inputs, plus move the active selection to match the diff's
# sheetView/selection element.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Computation")
$ws.Activate()

# Duty-cycle "on" time inputs for the two load rows (row 7 and row 11):
# 0.3 -> 1 (minutes on), which ripples through H/I/J on those rows.
$ws.Range("F7").Value = 1
$ws.Range("F11").Value = 1

# Battery / cycle count input used in the capacity-fade block.
$ws.Range("I18").Value = 3

# Efficiency input used at the bottom of the sheet.
$ws.Range("J36").Value = 0.8

# Match the author's final selection/scroll position from the diff
# (sheetView selection activeCell/sqref = J37; topLeftCell=A13 reflects
# the scrolled viewport, which follows from navigating/selecting there).
$ws.Range("A13").Select() | Out-Null
$ws.Range("J37").Select() | Out-Null
